# JobPlanning.xlsx – rotate the "Activity / ScheduledResource / PlannedQty /
# Plant / ProductionDivision" block among rows 2-4 and rows 7-9.
#
# Row 2 gets what row 4 had, row 3 gets what row 2 had, row 4 gets what row 3
# had (same cyclic rotation for rows 7/8/9).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cutterMulti = "406-45`" Polar 115ED Cutter`n404-45`" Polar 115EMC Cutter`n405-54`" Polar 137EMC Cutter`n402-45`" Polar 115EMC Cutter`n403-54`" Polar 137ED Cutter"

# --- Row 2 (becomes the "Cut" row) ---
$ws.Range("B2").Value = "Cut"
$ws.Range("D2").Value = "406-45`" Polar 115ED Cutter"
$ws.Range("G2").Value = "17"
$ws.Range("L2").Value = "406-45`" Polar 115ED Cutter"
$ws.Range("M2").Value = $cutterMulti

# --- Row 3 (becomes the "-" / Press Approval Task row) ---
$ws.Range("B3").Value = "-"
$ws.Range("D3").Value = "169-Press Approval Task "
$ws.Range("G3").Value = "24"
$ws.Range("L3").Value = "Press Approval Task"
$ws.Range("M3").Value = "169-Press Approval Task "

# --- Row 4 (becomes the "Digital Print F 4x0" row) ---
$ws.Range("B4").Value = "Digital Print F 4x0"
$ws.Range("D4").Value = "Versafire Heidelberg"
$ws.Range("G4").Value = "18"
$ws.Range("L4").Value = "Versafire Heidelberg"
$ws.Range("M4").Value = "Versafire Heidelberg"

# --- Row 7 (becomes the "-" / Press Approval Task row) ---
$ws.Range("B7").Value = "-"
$ws.Range("D7").Value = "169-Press Approval Task "
$ws.Range("G7").Value = "41"
$ws.Range("L7").Value = "Press Approval Task"
$ws.Range("M7").Value = "169-Press Approval Task "

# --- Row 8 (becomes the "Cut" row) ---
$ws.Range("B8").Value = "Cut"
$ws.Range("D8").Value = "406-45`" Polar 115ED Cutter"
$ws.Range("G8").Value = "34"
$ws.Range("L8").Value = "406-45`" Polar 115ED Cutter"
$ws.Range("M8").Value = $cutterMulti

# --- Row 9 (becomes the "Digital Print F 4x0" row) ---
$ws.Range("B9").Value = "Digital Print F 4x0"
$ws.Range("D9").Value = "Versafire Heidelberg"
$ws.Range("G9").Value = "35"
$ws.Range("L9").Value = "Versafire Heidelberg"
$ws.Range("M9").Value = "Versafire Heidelberg"
